# All 4 years scope emissions added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabeling ---
# Column C now holds "emissions" figures (previously "Target Max Electricity kWh per anum"),
# column D keeps the "Star Rating" label.
$ws.Range("C1").Value = "emissions"
$ws.Range("D1").Value = "Star Rating"

# --- Updated emissions values in column C (scope emissions for the 4 years) ---
$ws.Range("C2").Value = 44267.1
$ws.Range("C3").Value = 53915.8
$ws.Range("C4").Value = 63564.6
$ws.Range("C5").Value = 88534.3
$ws.Range("C6").Value = 107831.7
$ws.Range("C7").Value = 127129.1
$ws.Range("C8").Value = 132801.4
$ws.Range("C9").Value = 161747.5
$ws.Range("C10").Value = 190693.7
$ws.Range("C11").Value = 177068.5
$ws.Range("C12").Value = 215663.4
$ws.Range("C13").Value = 254258.3
$ws.Range("C14").Value = 221335.6
$ws.Range("C15").Value = 269579.2
$ws.Range("C16").Value = 317822.8

# --- Number format: column C now shows two decimal places ---
$ws.Range("C2:C16").NumberFormat = "#,##0.00"

# --- Column C widened to fit the new "emissions" values/header ---
# (target stored width is 19.33203125; COM ColumnWidth excludes the ~5/6 char
# cell-padding that Excel adds when persisting the <col> width attribute)
$ws.Columns.Item(3).ColumnWidth = 18.498697916666668

# --- Selection moved to E4 ---
$ws.Range("E4").Select()
